# RaFHY_pH9_NoScript.xlsx - "Fit FHY with new data points, data tweaking needed"
#
# Data edit: the counted solid mass (column O, 'Bottle Results') for rows
# 5-19 is corrected from 0.0225 g to 0.0226 g. For rows 5-7 (the "5 mg/L"
# triplicate) this pushes S (Cs) above U (total activity), so the
# fSorb/dCs column (W) goes negative there; the analyst respond by
# re-pointing Y5:Y7 (Cs used downstream) straight at S5:S7 instead of the
# shared "(U-Y*E)/Q" style estimate, matching the pattern already used by
# the Y8:Y19 block.
#
# Everything else (S, T, W, AA, AB on 'Bottle Results'; 'Sheet1';
# 'Averaged Results') is formula-driven off these inputs and recalculates
# automatically.

$wb = $excel.ActiveWorkbook
$bottle = $wb.Worksheets.Item("Bottle Results")

# --- Column O (Counted Solid Mass (g)) : 0.0225 -> 0.0226 for rows 5-19 ---
$bottle.Range("O5:O19").Value = 0.0226

# --- Y5:Y7 now read straight off S (instead of the shared W3:W19-style formula) ---
$bottle.Range("Y5").Formula = "=S5"
$bottle.Range("Y6").Formula = "=S6"
$bottle.Range("Y7").Formula = "=S7"

# --- View / selection state -------------------------------------------------

# 'Count->Actual Activity' : selection moved from C20 to H23
$countActivity = $wb.Worksheets.Item("Count->Actual Activity")
$countActivity.Activate()
$countActivity.Range("H23").Select()

# 'Sheet1' : selection moved from D11 to D7
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Activate()
$sheet1.Range("D7").Select()

# 'Averaged Results' : selection moved from B7 to B4:H7 (no longer the active tab)
$averaged = $wb.Worksheets.Item("Averaged Results")
$averaged.Activate()
$averaged.Range("B4:H7").Select()

# 'Bottle Results' : becomes the active tab; frozen-pane selection moves
# from Y21 to W20
$bottle.Activate()
$bottle.Range("W20").Select()
